$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Renumber the "Sl. No." column for existing rows 3-7 (1,1,2,2,2 -> 2,3,4,5,6) ---
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

# --- New rows for the "Keras ANN" classifier ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Keras ANN"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 20
$ws.Range("E8").Value = "without mfcc 14-20"
$ws.Range("F8").Value = 0.66

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Keras ANN"
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 35
$ws.Range("E9").Value = "rmse, mfcc mean&var(upto 13)"
$ws.Range("E9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 30
$ws.Range("F9").Value = 0.6733

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Keras ANN"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 35
$ws.Range("E10").Value = "rmse, mfcc mean&var(upto 13)"
$ws.Range("E10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 30
$ws.Range("F10").Value = 0.68

$ws.Range("A11").Value = 10

# --- Scroll / selection to match the final view ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("E10").Select()
